$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Special setups" header value from I1 down to I2, joining the
# rest of the header row (2)
$ws.Range("I1").ClearContents()
$ws.Range("I2").Value = "Special setups"

# Update the active selection to I3, as in the saved workbook
$ws.Range("I3").Select()

$wb.Save()
